# Insert a new row at position 57, shifting existing rows 57-82 down to 58-83,
# and populate the new row with this week's price data for
# Femacal de La Calera - Bruselas (repollito).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("57:57").Insert()

$row = 57
$ws.Cells.Item($row, 1).Value = 3
$ws.Cells.Item($row, 2).Value = "Femacal de La Calera"
$ws.Cells.Item($row, 3).Value = "Coquimbo"
$ws.Cells.Item($row, 4).Value = 44806
$ws.Cells.Item($row, 5).Value = 5
$ws.Cells.Item($row, 6).Value = 100112035
$ws.Cells.Item($row, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item($row, 8).Value = "Sin especificar"
$ws.Cells.Item($row, 9).Value = "Primera"
$ws.Cells.Item($row, 10).Value = 45
$ws.Cells.Item($row, 11).Value = 16000
$ws.Cells.Item($row, 12).Value = 16000
$ws.Cells.Item($row, 13).Value = 16000
$ws.Cells.Item($row, 14).Value = "`$/malla 15 kilos"
$ws.Cells.Item($row, 15).Value = "Provincia de Quillota"
$ws.Cells.Item($row, 16).Value = 1067
$ws.Cells.Item($row, 17).Value = 15
$ws.Cells.Item($row, 18).Value = "Hortaliza"
